$wb = $excel.ActiveWorkbook

# --- Capture worksheet references up front (by original name) ---
$wsPackages   = $wb.Worksheets.Item("packages")
$wsEntities   = $wb.Worksheets.Item("entities")
$wsAttributes = $wb.Worksheets.Item("attributes")
$wsAuthor     = $wb.Worksheets.Item("onetomany_author")
$wsBook       = $wb.Worksheets.Item("onetomany_book")
$wsNode       = $wb.Worksheets.Item("onetomany_node")

# --- Rename the one-to-many entity sheets to the new "it_emx_" prefixed names ---
$wsAuthor.Name = "it_emx_onetomany_author"
$wsBook.Name   = "it_emx_onetomany_book"
$wsNode.Name   = "it_emx_onetomany_node"

# --- packages sheet: rename "onetomany" -> "it_emx_onetomany" and add package hierarchy rows ---
$wsPackages.Cells.Replace("onetomany", "it_emx_onetomany", 1)
$wsPackages.Range("C2").Value = "it_emx"
$wsPackages.Range("A3").Value = "it_emx"
$wsPackages.Range("C3").Value = "it"
$wsPackages.Range("C3").Font.Color = 0
$wsPackages.Range("A4").Value = "it"
$wsPackages.Columns.Item(1).ColumnWidth = 15

# --- entities sheet: the "onetomany" package text also needs updating ---
$wsEntities.Cells.Replace("onetomany", "it_emx_onetomany", 1)
$wsEntities.Columns.Item(3).ColumnWidth = 15

# --- attributes sheet: rename entity references to the new sheet/entity names ---
$wsAttributes.Cells.Replace("onetomany_book", "it_emx_onetomany_book", 1)
$wsAttributes.Cells.Replace("onetomany_author", "it_emx_onetomany_author", 1)
$wsAttributes.Cells.Replace("onetomany_node", "it_emx_onetomany_node", 1)
$wsAttributes.Columns.Item(2).ColumnWidth = 21
$wsAttributes.Columns.Item(6).ColumnWidth = 21
$wsAttributes.Range("D18").Select()

# --- selections recorded by the author while editing ---
$wsPackages.Range("B8").Select()

# --- the node sheet becomes the active tab at the end ---
$wsNode.Activate()
$wsNode.Range("I22").Select()

Write-Host "edit complete"
